$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.537.33'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.08%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.857.30'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.65%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.005'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.16%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.40'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.03%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4642'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.13%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3885'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.32%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07887'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9711'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -0.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '22.18'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.22%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.800.22'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.16%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.723'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -1.42%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.925'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.22%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.06907'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.34%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '88.94'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.20%  '

$ws.Range("E17").Value = '  +0.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009999'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.02%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.82'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.06%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.003'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.10%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '28.553.39'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.04%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.323'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.07'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.82%  '

$ws.Range("E24").Value = '  -3.20%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.083.59'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.38%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '154.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +0.89%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.23%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.804'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.73%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.991'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.36%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '119.38'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.73%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09305'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.83%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.9400'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -2.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.302'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.01%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.333'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.56%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.326'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -3.95%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05848'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -3.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02124'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.87%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.148'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.825'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.58%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5633'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.940'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.63%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1771'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.09%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.07351'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +3.93%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '11.69'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.08%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5308'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.89%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.155'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -9.77%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.141'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -8.81%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.845'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.34%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '114.05'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.66%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.002'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.339'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.21%  '
